$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New/updated data for rows 2-9, columns A-T (20 columns)
# Column order: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = numeric metrics
$data = @(
    @("ECs", "Vip", "Adcyap1r1", "ECs", 2, 1, 1.38698, 2.77396, 0.5967131687073423, 0.4965808999056411, 1, 0.5, 0.100211, 0.200422, 0.01587220340208096, 0.01193595905534038, 0.13899065278, 0.5559626111199999, 0.009471152786423185, 0.005927169288937813),
    @("ECs", "Vip", "Adcyap1r1", "FAPs", 2, 1, 1.38698, 2.77396, 0.5967131687073423, 0.4965808999056411, 3, 1, 3.986394333333333, 11.959183, 0.6313963706536166, 0.7122188114245082, 5.529049212446666, 33.17429527468, 0.3767625290430351, 0.3536742583069084),
    @("ECs", "Vip", "Adcyap1r1", "Neutro", 2, 1, 1.38698, 2.77396, 0.5967131687073423, 0.4965808999056411, 1, 0.3333333333333333, 0.1778183333333333, 0.533455, 0.02816426096222669, 0.03176945164635921, 0.2466304719666666, 1.4797828318, 0.01680598540307079, 0.0157761028880578),
    @("ECs", "Vip", "Adcyap1r1", "sCs", 2, 1, 1.38698, 2.77396, 0.5967131687073423, 0.4965808999056411, 2, 1, 2.0491925, 4.098385, 0.3245671649820757, 0.2440757778737923, 2.84218901365, 11.3687560546, 0.1936735014748131, 0.1212033694217372),
    @("Neutro", "Vip", "Adcyap1r1", "ECs", 3, 1, 0.9373863333333334, 2.812159, 0.4032868312926577, 0.503419100094359, 1, 0.5, 0.100211, 0.200422, 0.01587220340208096, 0.01193595905534038, 0.09393642184966668, 0.563618531098, 0.006401050615657771, 0.00600878976640257),
    @("Neutro", "Vip", "Adcyap1r1", "FAPs", 3, 1, 0.9373863333333334, 2.812159, 0.4032868312926577, 0.503419100094359, 3, 1, 3.986394333333333, 11.959183, 0.6313963706536166, 0.7122188114245082, 3.736791567344111, 33.631124106097, 0.2546338416105815, 0.3585445531175999),
    @("Neutro", "Vip", "Adcyap1r1", "Neutro", 3, 1, 0.9373863333333334, 2.812159, 0.4032868312926577, 0.503419100094359, 1, 0.3333333333333333, 0.1778183333333333, 0.533455, 0.02816426096222669, 0.03176945164635921, 0.1666844754827778, 1.500160279345, 0.0113582755591559, 0.0159933487583014),
    @("Neutro", "Vip", "Adcyap1r1", "sCs", 3, 1, 0.9373863333333334, 2.812159, 0.4032868312926577, 0.503419100094359, 2, 1, 2.0491925, 4.098385, 0.3245671649820757, 0.2440757778737923, 1.920885043869167, 11.525310263215, 0.1308936635072626, 0.1228724084520552)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$j]
    }
}
